$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header to bold
$ws.Range("A1").Value = "Product"
$ws.Range("A1").Font.Bold = $true

# Update data rows
$ws.Range("A2").Value = "Business Loan"
$ws.Range("A3").Value = "Car Loan"
$ws.Range("A4").Value = "Home Loan"
$ws.Range("A5").Value = "Student Loan"

# Remove the now-unused Password column
$ws.Columns("B").Delete()

# Match the selection state saved in the workbook
$null = $ws.Range("A3").Select()

# Match page setup orientation
$ws.PageSetup.Orientation = 1
